$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert a new row above the old row 27 (this shifts the "Notes" block
# and everything below it down by one row, matching the target layout).
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the "partial first year" helper
# inputs: a starting month label and the fraction of the year remaining.
$ws.Cells.Item(27, 1).Value2 = "Starting month"
$ws.Cells.Item(27, 2).Value2 = "April"
$ws.Cells.Item(27, 3).Value2 = 0.75

# Scale the first-year (partial year) figures by the new C27 fraction, and
# adjust the dependent "de-annualizing" formulas so later years are
# unaffected by the partial first year.
$ws.Range("B3").Formula = "=70*1500*C27"
$ws.Range("C3").Formula = "=B3/C27"

$ws.Range("B13").Formula = "=ROUND((9500 + 0.25*((B12/C27)- 88000) + B3 * 0.0765),0)"

$ws.Range("B15").Formula = "=72000*C27"
$ws.Range("C15").Formula = "=ROUND((B15/C27)*1.03,0)"

$ws.Range("B16").Formula = "=12000*C27"
$ws.Range("C16").Formula = "=ROUND((B16/C27)*1.08,0)"

# Recalculate everything so cached values stay consistent with the formulas.
$excel.CalculateFullRebuild()

# Leave the selection on B15, mirroring the cell the edit focused on.
$ws.Range("B15").Select()
